$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 21, shifting the existing rows 21-105 down to 22-106.
$ws.Rows("21:21").Insert()

# Populate the newly inserted row 21 with this week's data point.
# Most columns carry the same constant/contextual values as their neighbours;
# only the date (D), min/max/avg price (K/L/M) and $/Kg price (P) are new.
$ws.Range("A21").Value = 6
$ws.Range("B21").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C21").Value = "Metropolitana"
$ws.Range("D21").Value = 44453
$ws.Range("E21").Value = 13
$ws.Range("F21").Value = 100112029
$ws.Range("G21").Value = "Orégano"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 32
$ws.Range("K21").Value = 8000
$ws.Range("L21").Value = 9000
$ws.Range("M21").Value = 8438
$ws.Range("N21").Value = "$/docena de atados"
$ws.Range("O21").Value = "Región Metropolitana"
$ws.Range("P21").Value = 2813
$ws.Range("Q21").Value = 3
$ws.Range("R21").Value = "Hortaliza"
